# MiMalla.xlsx edit script
# - bump the curricular-year header (B1) 2018 -> 2010
# - small prerequisite/semester corrections (clique change)
# - clear the placeholder rows 23-29 (front-end will populate these from the
#   "get available sections" function instead of hard-coded data)
# - move the active selection to A3:F22 (the now-complete data block)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header year
$ws.Range("B1").Value = "2010"

# Targeted value corrections
$ws.Range("E4").Value = "7, 8, 33"
$ws.Range("E5").Value = "53"
$ws.Range("E7").Value = "38,21"
$ws.Range("B12").Value = "CFG-1"
$ws.Range("E12").Value = "53"
$ws.Range("E15").Value = "53"
$ws.Range("E19").Value = "53"
$ws.Range("B22").Value = "CFG-2"
$ws.Range("E22").Value = "53"

# Clear out the old hard-coded rows 23-29 (keep their formatting)
$ws.Range("A23:F29").ClearContents()

# Update selection to the active data block
$ws.Range("A3:F22").Select()
